$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo: missing space before "TAŞ"
$ws.Range("A88").Value = "ŞİŞME YELEK 5010 YELEK TAŞ"

# Drop the "DERİ" word from these two product names
$ws.Range("A102").Value = "ŞİŞME YELEK   3004 ALARA  YELEK SİYAH"
$ws.Range("A103").Value = "ŞİŞME YELEK   3004 ALARA  YELEK HAKİ"

# Append two new products (rows 104-105), column by column
$ws.Range("A104").Value = "ERKEK ŞİŞME YELEK 5431 HAKİ"
$ws.Range("A105").Value = "ERKEK ŞİŞME YELEK 5431 SİYAH"

$ws.Range("B104").Value = "575 TL"
$ws.Range("B105").Value = "575 TL"

$ws.Range("C104").Value = "Yelek"
$ws.Range("C105").Value = "Yelek"

$ws.Range("D104").Value = "erkek yelek5431.jpg"
$ws.Range("D105").Value = "erkekyeleksiyah.jpg"

$ws.Range("E104").Value = "S-M-L-XL-2XL-3XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("E105").Value = "S-M-L-XL-2XL-3XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

$ws.Range("F104").Value = "Var"
$ws.Range("F105").Value = "Var"

$ws.Range("E108").Select()
